$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new item row ("HIBIOTIC 1GM 16 TAB") right before "LAGAR 15 ML
# DROPS" (the row that was row 15), pushing that row and everything below it
# down by one row.
# ---------------------------------------------------------------------------
$ws.Range("A15:N15").Insert(-4121)  # xlShiftDown

# Copy the cell formatting (styles/borders/fill) from the row just below
# (the row that now holds the old row-15 data) onto the freshly inserted,
# still-blank row 15 so it looks like every other data row.
$ws.Range("A16:G16").Copy()
$ws.Range("A15:G15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H16:K16").Copy()
$ws.Range("H15:K15").PasteSpecial(-4122)
$ws.Range("L16:M16").Copy()
$ws.Range("L15:M15").PasteSpecial(-4122)
$ws.Range("N16").Copy()
$ws.Range("N15").PasteSpecial(-4122)

# Re-create the merged cells for the new row, matching the layout used by
# every other data row (B:G, H:K, L:M).
$ws.Range("B15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()

# Populate the new row's values.
$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "HIBIOTIC 1GM 16 TAB"
$ws.Range("H15").Value = "1:1"
$ws.Range("L15").Value = 86.5
$ws.Range("N15").Value = 0.5

# The row-insert also shifted the column-A running counter ("م") down along
# with the rest of the row content; restore it to the plain 1..31 sequence
# expected for rows 4-34 (row 16 -> 13, row 17 -> 14, ... row 34 -> 31).
for ($r = 16; $r -le 34; $r++) {
    $ws.Range("A$r").Value = $r - 3
}

# Update the grand-total cell (now on row 35) to include the new row's price.
$ws.Range("K35").Value = 2298.2199999999998

# Restore the per-row heights so they again match the (content-driven)
# heights used by the row that used to occupy that slot - i.e. every row
# from 16 downward keeps the height that the row above it had before the
# insert, row 15 keeps its original height, and the totals row grows
# slightly to fit its new position.
$ws.Rows(15).RowHeight = 25.5
$ws.Rows(16).RowHeight = 24.75
$ws.Rows(17).RowHeight = 25.5
$ws.Rows(18).RowHeight = 25.5
$ws.Rows(19).RowHeight = 24.75
$ws.Rows(20).RowHeight = 25.5
$ws.Rows(21).RowHeight = 24.75
$ws.Rows(22).RowHeight = 25.5
$ws.Rows(23).RowHeight = 25.5
$ws.Rows(24).RowHeight = 24.75
$ws.Rows(25).RowHeight = 25.5
$ws.Rows(26).RowHeight = 24.75
$ws.Rows(27).RowHeight = 25.5
$ws.Rows(28).RowHeight = 25.5
$ws.Rows(29).RowHeight = 24.75
$ws.Rows(30).RowHeight = 25.5
$ws.Rows(31).RowHeight = 24.75
$ws.Rows(32).RowHeight = 25.5
$ws.Rows(33).RowHeight = 25.5
$ws.Rows(34).RowHeight = 24.75
$ws.Rows(35).RowHeight = 26.25
$ws.Rows(36).RowHeight = 16.5

Write-Host "Edit applied successfully"
